# Update data: 4 December 2020
# Adds the newest observation (date serial 44136 = 2020-11-01) to both the
# "Canada" sheet (one row) and the "Province" sheet (ten rows, one per
# province) and updates the sheet views/selections to point at the new data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Canada": append row 12
# ---------------------------------------------------------------------
$wsCanada = $wb.Worksheets.Item("Canada")

# Clone the formatting of the last existing data row (11) down onto the new
# row 12 so the date/style formatting (style index carrying the date number
# format) matches the rest of the column without minting a new style.
$wsCanada.Range("A11:B11").Copy()
$wsCanada.Range("A12:B12").PasteSpecial(-4122)

$wsCanada.Range("A12").Value = 44136
$wsCanada.Range("B12").Value = "Canada"
$wsCanada.Range("C12").Value = 45.6
$wsCanada.Range("D12").Value = 1735.2

# ---------------------------------------------------------------------
# Sheet "Province": append rows 102-111 (one per province, same order as
# every prior monthly block)
# ---------------------------------------------------------------------
$wsProvince = $wb.Worksheets.Item("Province")

# Clone the formatting of the previous month's 10-row block (rows 92-101)
# down onto the new block (rows 102-111).
$wsProvince.Range("A92:D101").Copy()
$wsProvince.Range("A102:D111").PasteSpecial(-4122)

$wsProvince.Range("A102").Value = 44136
$wsProvince.Range("B102").Value = "Newfoundland & Labrador"
$wsProvince.Range("C102").Value = 8.7
$wsProvince.Range("D102").Value = 31.3

$wsProvince.Range("A103").Value = 44136
$wsProvince.Range("B103").Value = "Prince Edward Island"
$wsProvince.Range("C103").Value = 27.5
$wsProvince.Range("D103").Value = 8.8

$wsProvince.Range("A104").Value = 44136
$wsProvince.Range("B104").Value = "Nova Scotia"
$wsProvince.Range("C104").Value = -19.1
$wsProvince.Range("D104").Value = 32.2

$wsProvince.Range("A105").Value = 44136
$wsProvince.Range("B105").Value = "New Brunswick"
$wsProvince.Range("C105").Value = 24.4
$wsProvince.Range("D105").Value = 37.7

$wsProvince.Range("A106").Value = 44136
$wsProvince.Range("B106").Value = "Quebec"
$wsProvince.Range("C106").Value = 29.4
$wsProvince.Range("D106").Value = 329.7

$wsProvince.Range("A107").Value = 44136
$wsProvince.Range("B107").Value = "Ontario"
$wsProvince.Range("C107").Value = 66.6
$wsProvince.Range("D107").Value = 733.7

$wsProvince.Range("A108").Value = 44136
$wsProvince.Range("B108").Value = "Manitoba"
$wsProvince.Range("C108").Value = 32.6
$wsProvince.Range("D108").Value = 50.9

$wsProvince.Range("A109").Value = 44136
$wsProvince.Range("B109").Value = "Saskatchewan"
$wsProvince.Range("C109").Value = 16.5
$wsProvince.Range("D109").Value = 41.7

$wsProvince.Range("A110").Value = 44136
$wsProvince.Range("B110").Value = "Alberta"
$wsProvince.Range("C110").Value = 51.6
$wsProvince.Range("D110").Value = 276.7

$wsProvince.Range("A111").Value = 44136
$wsProvince.Range("B111").Value = "British Columbia"
$wsProvince.Range("C111").Value = 43
$wsProvince.Range("D111").Value = 192.4

# ---------------------------------------------------------------------
# Sheet views / selections: point them at the newly appended rows. Select
# the "Canada" sheet's new row first so that "Province" ends up being the
# final active / tab-selected sheet, matching the original workbook state.
# ---------------------------------------------------------------------
$wsCanada.Range("A12:D12").Select()
$wsProvince.Range("A102:D111").Select()

# Scroll the "Province" view down so the freshly appended block is visible
# (mirrors the author scrolling to A89 before saving).
$excel.ActiveWindow.ScrollRow = 89
$excel.ActiveWindow.ScrollColumn = 1
